# excelcy/tests/data/test_data_28.xlsx
#  - "Add error codes" / "Fix data storage issues" / "Improve tests"
#
# The workbook has 3 sheets: source (1), prepare (2), config (3).
#
# 1) source!C1 and prepare!C1 header renamed from "text" -> "value"
# 2) prepare sheet had a stray/stale selection left over at C3; clear it
# 3) config sheet: the "nlp_name" row (row 3) is removed entirely, and
#    the final row's value (train_autosave) is flipped from TRUE to FALSE

$wb = $excel.ActiveWorkbook

# --- sheet "source": rename column header text -> value ---
$wsSource = $wb.Worksheets.Item("source")
$wsSource.Range("C1").Value = "value"

# --- sheet "prepare": rename column header text -> value ---
$wsPrepare = $wb.Worksheets.Item("prepare")
$wsPrepare.Range("C1").Value = "value"
# clear the stale selection that pointed at C3
$wsPrepare.Range("A1").Select()

# --- sheet "config": drop the "nlp_name" row, fix train_autosave value ---
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Rows.Item(3).Delete()
$wsConfig.Cells.Item(7, 2).Value = $false

# keep "config" as the active/visible tab (matches original workbook state)
$wsConfig.Activate()
